$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4249093333333334
$ws.Range("H2").Value = 1.274728
$ws.Range("I2").Value = 0.06780552842016908
$ws.Range("J2").Value = 0.06780552842016908
$ws.Range("M2").Value = 192.036265
$ws.Range("N2").Value = 576.108795
$ws.Range("O2").Value = 0.9537264487607444
$ws.Range("P2").Value = 0.9537264487607444
$ws.Range("Q2").Value = 81.59800133697333
$ws.Range("R2").Value = 734.38201203276
$ws.Range("S2").Value = 0.06466792582651359
$ws.Range("T2").Value = 0.06466792582651359

# Row 3
$ws.Range("G3").Value = 0.4249093333333334
$ws.Range("H3").Value = 1.274728
$ws.Range("I3").Value = 0.06780552842016908
$ws.Range("J3").Value = 0.06780552842016908
$ws.Range("O3").Value = 0.005072929450888834
$ws.Range("P3").Value = 0.005072929450888834
$ws.Range("Q3").Value = 0.4340247716248889
$ws.Range("R3").Value = 3.906222944624001
$ws.Range("S3").Value = 0.0003439726620557555
$ws.Range("T3").Value = 0.0003439726620557555

# Row 4
$ws.Range("G4").Value = 0.4249093333333334
$ws.Range("H4").Value = 1.274728
$ws.Range("I4").Value = 0.06780552842016908
$ws.Range("J4").Value = 0.06780552842016908
$ws.Range("M4").Value = 8.295893999999999
$ws.Range("N4").Value = 24.887682
$ws.Range("O4").Value = 0.04120062178836673
$ws.Range("P4").Value = 0.04120062178836673
$ws.Range("Q4").Value = 3.525002788944
$ws.Range("R4").Value = 31.725025100496
$ws.Range("S4").Value = 0.002793629931599738
$ws.Range("T4").Value = 0.002793629931599738

# Row 5
$ws.Range("I5").Value = 0.4624930683973976
$ws.Range("J5").Value = 0.4624930683973975
$ws.Range("M5").Value = 192.036265
$ws.Range("N5").Value = 576.108795
$ws.Range("O5").Value = 0.9537264487607444
$ws.Range("P5").Value = 0.9537264487607444
$ws.Range("Q5").Value = 556.5698091692217
$ws.Range("R5").Value = 5009.128282522995
$ws.Range("S5").Value = 0.4410918716991101
$ws.Range("T5").Value = 0.44109187169911

# Row 6
$ws.Range("I6").Value = 0.4624930683973976
$ws.Range("J6").Value = 0.4624930683973975
$ws.Range("O6").Value = 0.005072929450888834
$ws.Range("P6").Value = 0.005072929450888834
$ws.Range("S6").Value = 0.002346194707505102
$ws.Range("T6").Value = 0.002346194707505102

# Row 7
$ws.Range("I7").Value = 0.4624930683973976
$ws.Range("J7").Value = 0.4624930683973975
$ws.Range("M7").Value = 8.295893999999999
$ws.Range("N7").Value = 24.887682
$ws.Range("O7").Value = 0.04120062178836673
$ws.Range("P7").Value = 0.04120062178836673
$ws.Range("Q7").Value = 24.043605203778
$ws.Range("R7").Value = 216.392446834002
$ws.Range("S7").Value = 0.0190550019907824
$ws.Range("T7").Value = 0.0190550019907824

# Row 8
$ws.Range("G8").Value = 2.943425333333333
$ws.Range("H8").Value = 8.830276
$ws.Range("I8").Value = 0.4697014031824334
$ws.Range("J8").Value = 0.4697014031824334
$ws.Range("M8").Value = 192.036265
$ws.Range("N8").Value = 576.108795
$ws.Range("O8").Value = 0.9537264487607444
$ws.Range("P8").Value = 0.9537264487607444
$ws.Range("Q8").Value = 565.2444073197133
$ws.Range("R8").Value = 5087.19966587742
$ws.Range("S8").Value = 0.4479666512351209
$ws.Range("T8").Value = 0.4479666512351208

# Row 9
$ws.Range("G9").Value = 2.943425333333333
$ws.Range("H9").Value = 8.830276
$ws.Range("I9").Value = 0.4697014031824334
$ws.Range("J9").Value = 0.4697014031824334
$ws.Range("O9").Value = 0.005072929450888834
$ws.Range("P9").Value = 0.005072929450888834
$ws.Range("Q9").Value = 3.006569655867556
$ws.Range("R9").Value = 27.059126902808
$ws.Range("S9").Value = 0.002382762081327977
$ws.Range("T9").Value = 0.002382762081327976

# Row 10
$ws.Range("G10").Value = 2.943425333333333
$ws.Range("H10").Value = 8.830276
$ws.Range("I10").Value = 0.4697014031824334
$ws.Range("J10").Value = 0.4697014031824334
$ws.Range("M10").Value = 8.295893999999999
$ws.Range("N10").Value = 24.887682
$ws.Range("O10").Value = 0.04120062178836673
$ws.Range("P10").Value = 0.04120062178836673
$ws.Range("Q10").Value = 24.418344562248
$ws.Range("R10").Value = 219.765101060232
$ws.Range("S10").Value = 0.01935198986598459
$ws.Range("T10").Value = 0.01935198986598459
